# Lower3 worksheet update: replace existing data rows (2-7) with new values
# and append additional rows (8-29) of Time / Lower3 / MA data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(44090.02083333334, 2.5693319781592, $null)
    ,@(44090.66666666666, $null, 2.626344181182382)
    ,@(44105.67708333334, 2.767587059676764, $null)
    ,@(44105.9375, $null, 2.830988566391733)
    ,@(44106.42708333334, 2.46493470209716, $null)
    ,@(44106.58333333334, $null, 2.630449014213091)
    ,@(44110.46875, 2.388156873432133, $null)
    ,@(44111.21875, $null, 2.221048679547392)
    ,@(44124.33333333334, 1.933802760926178, $null)
    ,@(44124.84375, $null, 1.915600803314232)
    ,@(44130.6875, 1.661015460463576, $null)
    ,@(44131.08333333334, $null, 1.708347335846987)
    ,@(44161.14583333334, 1.809724619994121, $null)
    ,@(44161.85416666666, $null, 1.819467800485715)
    ,@(44173.39583333334, 1.681061761354724, $null)
    ,@(44173.67708333334, $null, 1.731193274007974)
    ,@(44188.45833333334, 1.285333792524691, $null)
    ,@(44189.0625, $null, 1.265586180957891)
    ,@(44218.04166666666, 2.659821852105991, $null)
    ,@(44218.15625, $null, 3.101020120276233)
    ,@(44242.09375, 7.801765455997229, $null)
    ,@(44242.35416666666, $null, 8.369873372139002)
    ,@(44304.13541666666, 21.92602548626203, $null)
    ,@(44304.25, $null, 24.35933694351101)
    ,@(44326.83333333334, 39.24104286753596, $null)
    ,@(44326.90625, $null, 43.53030316279118)
    ,@(44335.53125, 29.68888367864963, $null)
    ,@(44335.63541666666, $null, 45.68648690290231)
)

# Row 2 already carries the date number format / style we want to reuse
# for every Time cell in column A (including the newly appended rows).
$firstDateCell = $ws.Range("A2")

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    if ($row[2] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $r = $r + 1
}

$lastRow = $r - 1

# Apply the existing date formatting (style) used in A2:A7 to the newly
# created cells A8:A29, without disturbing the values just written.
$firstDateCell.Copy()
$ws.Range("A8:A" + $lastRow).PasteSpecial(-4122)
